# Change class labels in tukey input
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers J1 and K1, and add two new headers L1, M1
$ws.Range("J1").Value = "Peak_Fre_10"
$ws.Range("K1").Value = "Peak_Fre_20"
$ws.Range("L1").Value = "Peak_Dur_2"
$ws.Range("M1").Value = "Peak_Tim_2"

# Update selection to M7 to mirror the saved selection state in the target file
$ws.Range("M7").Select()
